# Refined metadata to be additional tab
#
# 1) Refresh the "panel_query_time" / time_taken stamps on the "data" sheet
#    (column F) to the new query run.
# 2) Append a new "metadata" worksheet (after "data") summarising the panel
#    query itself.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1) Update column F ("time_taken") timestamps on the data sheet -------
$newTimes = @(
    "2021-10-05 14:20:10.681329",
    "2021-10-05 14:20:10.681336",
    "2021-10-05 14:20:10.681339",
    "2021-10-05 14:20:10.681342",
    "2021-10-05 14:20:10.681345",
    "2021-10-05 14:20:10.681348",
    "2021-10-05 14:20:10.681350",
    "2021-10-05 14:20:10.681352",
    "2021-10-05 14:20:10.681355",
    "2021-10-05 14:20:10.681358",
    "2021-10-05 14:20:10.681360",
    "2021-10-05 14:20:10.681363",
    "2021-10-05 14:20:10.681365",
    "2021-10-05 14:20:10.681368",
    "2021-10-05 14:20:10.681370",
    "2021-10-05 14:20:10.681373",
    "2021-10-05 14:20:10.681375",
    "2021-10-05 14:20:10.681378",
    "2021-10-05 14:20:10.681380",
    "2021-10-05 14:20:10.681383",
    "2021-10-05 14:20:10.681385",
    "2021-10-05 14:20:10.681388",
    "2021-10-05 14:20:10.681390",
    "2021-10-05 14:20:10.681393",
    "2021-10-05 14:20:10.681396",
    "2021-10-05 14:20:10.681398"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $data.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2) Add the "metadata" worksheet, placed after "data" -----------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$meta = $wb.Worksheets.Add($null, $afterSheet)
$meta.Name = "metadata"

# Header row + the numeric index cell (A2) reuse the same bold/bordered
# style as the "data" sheet's header row -- copy it across instead of
# re-deriving the formatting (Font/Borders/Alignment) by hand, so the
# workbook keeps a single shared style instead of growing a near-duplicate.
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$data.Range("A2").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# Data row
$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 2).Value = "Familial breast cancer"
$meta.Cells.Item(2, 3).Value = 158

# "data_version" (1.14) is stored as text in the source data, not a number
# -- format the cell as Text first so Excel doesn't coerce it to a float.
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.14"

$meta.Cells.Item(2, 5).Value = "2021-07-15T09:16:24.340600Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:20:10.677537"
$meta.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/158/?format=json"
